# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to currentAveragePrice / Leve profit columns
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 912712.2
$ws.Range("I112").Value = 4098
$ws.Range("K112").Value = 12294
$ws.Range("M112").Value = -11186
$ws.Range("H137").Value = 1851.3334
$ws.Range("I137").Value = 895.8
$ws.Range("J137").Value = 2102.7896
$ws.Range("K137").Value = 2687.4
$ws.Range("L137").Value = 6308.3688
$ws.Range("M137").Value = -137.3999999999996
$ws.Range("N137").Value = -11408.3688
$ws.Range("H138").Value = 1669603.5
$ws.Range("I138").Value = 2155.5
$ws.Range("J138").Value = 5004499.5
$ws.Range("K138").Value = 6466.5
$ws.Range("L138").Value = 15013498.5
$ws.Range("M138").Value = -1326.5
$ws.Range("N138").Value = -15023778.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2951.3333
$ws.Range("I32").Value = 2932.7368
$ws.Range("K32").Value = 2932.7368
$ws.Range("M32").Value = -2645.7368
$ws.Range("H61").Value = 6397.409
$ws.Range("I61").Value = 9437
$ws.Range("J61").Value = 2749.9
$ws.Range("K61").Value = 9437
$ws.Range("L61").Value = 2749.9
$ws.Range("M61").Value = -9225
$ws.Range("N61").Value = -3173.9
$ws.Range("H74").Value = 5019.5
$ws.Range("I74").Value = 4799
$ws.Range("J74").Value = 5240
$ws.Range("K74").Value = 4799
$ws.Range("L74").Value = 5240
$ws.Range("M74").Value = -3925
$ws.Range("N74").Value = -6988
$ws.Range("H77").Value = 5019.5
$ws.Range("I77").Value = 4799
$ws.Range("J77").Value = 5240
$ws.Range("K77").Value = 23995
$ws.Range("L77").Value = 26200
$ws.Range("M77").Value = -19627
$ws.Range("N77").Value = -34936
$ws.Range("H122").Value = 2064.3157
$ws.Range("I122").Value = 2113.0312
$ws.Range("J122").Value = 1804.5
$ws.Range("K122").Value = 6339.0936
$ws.Range("L122").Value = 5413.5
$ws.Range("M122").Value = -3889.0936
$ws.Range("N122").Value = -10313.5
$ws.Range("H132").Value = 3440.0576
$ws.Range("I132").Value = 2235.8064
$ws.Range("K132").Value = 6707.4192
$ws.Range("M132").Value = -4177.4192
$ws.Range("H136").Value = 6397.409
$ws.Range("I136").Value = 9437
$ws.Range("J136").Value = 2749.9
$ws.Range("K136").Value = 28311
$ws.Range("L136").Value = 8249.700000000001
$ws.Range("M136").Value = -25761
$ws.Range("N136").Value = -13349.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4529.1113
$ws.Range("I107").Value = 5800
$ws.Range("J107").Value = 1987.3334
$ws.Range("K107").Value = 5800
$ws.Range("L107").Value = 1987.3334
$ws.Range("M107").Value = -3880
$ws.Range("N107").Value = -5827.3334
$ws.Range("H134").Value = 5532.8823
$ws.Range("I134").Value = 4216.72
$ws.Range("J134").Value = 9188.888999999999
$ws.Range("K134").Value = 12650.16
$ws.Range("L134").Value = 27566.667
$ws.Range("M134").Value = -10115.16
$ws.Range("N134").Value = -32636.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5841.1626
$ws.Range("I31").Value = 4690.5454
$ws.Range("K31").Value = 4690.5454
$ws.Range("M31").Value = -4395.5454
$ws.Range("H34").Value = 5841.1626
$ws.Range("I34").Value = 4690.5454
$ws.Range("K34").Value = 4690.5454
$ws.Range("M34").Value = -4488.5454
$ws.Range("H58").Value = 2942.5
$ws.Range("I58").Value = 4513.3335
$ws.Range("K58").Value = 4513.3335
$ws.Range("M58").Value = -4310.3335
$ws.Range("H132").Value = 1903.2307
$ws.Range("I132").Value = 1895.1666
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 5685.4998
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -3155.4998
$ws.Range("N132").Value = -11060
$ws.Range("H134").Value = 2035
$ws.Range("I134").Value = 2249
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 6747
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -4212
$ws.Range("N134").Value = -9570
$ws.Range("H136").Value = 2942.5
$ws.Range("I136").Value = 4513.3335
$ws.Range("K136").Value = 13540.0005
$ws.Range("M136").Value = -10990.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 958107.2
$ws.Range("I56").Value = 958107.2
$ws.Range("K56").Value = 958107.2
$ws.Range("M56").Value = -957577.2
$ws.Range("H68").Value = 7356.7144
$ws.Range("I68").Value = 1999
$ws.Range("J68").Value = 8249.666999999999
$ws.Range("K68").Value = 5997
$ws.Range("L68").Value = 24749.001
$ws.Range("M68").Value = -5186
$ws.Range("N68").Value = -26371.001
$ws.Range("H71").Value = 7356.7144
$ws.Range("I71").Value = 1999
$ws.Range("J71").Value = 8249.666999999999
$ws.Range("K71").Value = 17991
$ws.Range("L71").Value = 74247.003
$ws.Range("M71").Value = -13935
$ws.Range("N71").Value = -82359.003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 800.1389
$ws.Range("I97").Value = 508
$ws.Range("J97").Value = 1384.4166
$ws.Range("K97").Value = 508
$ws.Range("L97").Value = 1384.4166
$ws.Range("M97").Value = -12
$ws.Range("N97").Value = -2376.4166
$ws.Range("H102").Value = 28815.475
$ws.Range("I102").Value = 1617.6666
$ws.Range("J102").Value = 53293.5
$ws.Range("K102").Value = 1617.6666
$ws.Range("L102").Value = 53293.5
$ws.Range("M102").Value = 4.333399999999983
$ws.Range("N102").Value = -56537.5
$ws.Range("H113").Value = 1920.2593
$ws.Range("I113").Value = 1633.5333
$ws.Range("J113").Value = 2278.6667
$ws.Range("K113").Value = 1633.5333
$ws.Range("L113").Value = 2278.6667
$ws.Range("M113").Value = 536.4666999999999
$ws.Range("N113").Value = -6618.6667
$ws.Range("H132").Value = 3945.862
$ws.Range("I132").Value = 4377.2915
$ws.Range("J132").Value = 1875
$ws.Range("K132").Value = 13131.8745
$ws.Range("L132").Value = 5625
$ws.Range("M132").Value = -10601.8745
$ws.Range("N132").Value = -10685

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4623.537
$ws.Range("I132").Value = 4407.026
$ws.Range("J132").Value = 5186.467
$ws.Range("K132").Value = 13221.078
$ws.Range("L132").Value = 15559.401
$ws.Range("M132").Value = -10691.078
$ws.Range("N132").Value = -20619.401
$ws.Range("H136").Value = 3995.7307
$ws.Range("I136").Value = 3756.6191
$ws.Range("K136").Value = 11269.8573
$ws.Range("M136").Value = -8719.8573

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 149999
$ws.Range("J46").Value = 149999
$ws.Range("L46").Value = 149999
$ws.Range("N46").Value = -150461
$ws.Range("H113").Value = 3205920
$ws.Range("J113").Value = 1290
$ws.Range("L113").Value = 3870
$ws.Range("N113").Value = -8210
$ws.Range("H132").Value = 1712.1923
$ws.Range("I132").Value = 1550.8667
$ws.Range("K132").Value = 4652.6001
$ws.Range("M132").Value = -2122.6001
$ws.Range("H134").Value = 149999
$ws.Range("J134").Value = 149999
$ws.Range("L134").Value = 449997
$ws.Range("N134").Value = -455067
